# Add two new Saudi cities (Al Asyah and Al Hareeq) to the bottom of the
# cities table on Sheet1, mirroring the format/formulas of the existing
# rows (152 rows of data -> 154 rows of data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (borders/style) of the last existing data row (152)
# onto the two new rows before filling in values, so the new cells pick up
# the same bordered style used throughout the table instead of the default
# "no style" formatting.
$ws.Range("A152:G152").Copy()
$ws.Range("A153:G154").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 153: Al Asyah / الأسياح (Eastern Region / شرق المملكة)
$ws.Range("A153").Value = "Al Asyah"
$ws.Range("B153").Value = "Al Asyah"
$ws.Range("C153").Value = "الأسياح"
$ws.Range("D153").Value = 27.378982000000001
$ws.Range("E153").Value = 44.237861000000002
$ws.Range("F153").Value = "المنطقة الشرقية"
$ws.Range("G153").Value = "شرق المملكة"

# Row 154: Al Hareeq / الحريق (Riyadh Region / وسط المملكة)
$ws.Range("A154").Value = "Al Hareeq"
$ws.Range("B154").Value = "Al Hareeq"
$ws.Range("C154").Value = "الحريق"
$ws.Range("D154").Value = 23.624414000000002
$ws.Range("E154").Value = 46.511069999999997
$ws.Range("F154").Value = "منطقة الرياض"
$ws.Range("G154").Value = "وسط المملكة"

# Keep the sheet's selection consistent with the now-larger used range.
$ws.Range("A1:G154").Select()
